# "to cuda0 in sqlite" - update the GPU column on the rtsp camera list
# (the tab literally labelled "Sheet2" holds the camera table - the tab
# labelled "Sheet1" is a different, unrelated blank sheet): move cams
# 15/16/20/8/10/13/14/17/18/19 (rows 7-16) onto gpu 1, and cams 2/3
# (rows 2-3) back onto gpu 0. Also move the active selection to F6.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# Rows 2-3: gpu -> 0
$ws.Range("D2").Value = 0
$ws.Range("D3").Value = 0

# Rows 7-16: gpu -> 1
$ws.Range("D7:D16").Value = 1

# Match the saved selection in the file
$ws.Activate()
$ws.Range("F6").Select()
